# Add a new "Evaluated" (Yes/No) column to the "Others" sheet.
#
# Source row layout (before):
#   A: SAP ID | B..I: MCQ Q.1..7, Code Q.1 | J: Total
#   row2 -> SAP ID 001, row3 -> SAP ID 002
#
# Target: add column K "Evaluated" with "Yes" for SAP ID 001 and
# "No" for SAP ID 002, matching the column width (15) used by the
# other data columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Evaluated" column.
$ws.Range("K1").Value = "Evaluated"

# New values for each existing data row.
$ws.Range("K2").Value = "Yes"
$ws.Range("K3").Value = "No"

# Match the 15-character width used by the other data columns (B:J).
# ColumnWidth uses Excel's "characters" unit, which is offset from the
# raw stored width by the default font padding (~0.8333 chars), so use
# 15 - 5/6 to land on a stored width of exactly 15.
$ws.Columns("K").ColumnWidth = 14.166666666666666
